$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.875.93"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "2.448.54"
$ws.Range("E3").Value = "  -2.88%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").Value = "2.451.54"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.65%  "
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").Value = "2.886.10"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").Value = "57.795.07"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "2.448.49"
$ws.Range("E18").Value = "  -2.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "316.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.402"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.86%  "
$ws.Range("D30").Value = "0.0₃0737"
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("E33").Value = "  -5.04%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.94%  "
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.797"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "263.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.586"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.51%  "
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0212"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.47%  "
